# finsihed variant comparison pathway enrichment graph
#
# The EnrichNet results table (A1:E21) gets sorted by the "Hits" column
# (column B) in descending order - the user re-ran Data > Sort on the
# pathway-enrichment table. As a result of doing this interactively on
# the EnrichNet sheet, that sheet becomes the active sheet/tab (moving
# the active tab away from "Mummichog Output"), and the selection on
# EnrichNet moves to A2 (top-left of the now-reordered data).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("EnrichNet")

# Make EnrichNet the active sheet (this is what the user was working on).
$ws1.Activate()

# Sort the data rows (A2:E21) by Hits (column B), descending - headers in
# row 1 stay put.
$sortRange = $ws1.Range("A2:E21")
$sortKey = $ws1.Range("B2:B21")
$sortRange.Sort($sortKey, 2, $null, $null, 1, 0, 0, 0)

# Leave the selection on the first data cell, like after an interactive sort.
$ws1.Range("A2").Select()
